$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 8 (ano 2025) with refreshed customer metrics
$ws.Range("C8").Value = 868
$ws.Range("D8").Value = 146
$ws.Range("E8").Value = 722
$ws.Range("F8").Value = 5.988515176374077
$ws.Range("G8").Value = 83.17972350230414
$ws.Range("H8").Value = 16.82027649769585
